$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1 -and $parts[$parts.Count - 1] -eq "System") {
            $rest = $parts[0..($parts.Count - 2)]
            $newParts = @("System") + $rest
            $cell.Value2 = $newParts -join ", "
        }
    }
}
